# Generate Report for Archive
#
# 1) Status text update: "Ready for handoff" -> "In Translation"
#    This shared string is used by the Status-like column on every sheet:
#      - Overview sheet: columns E (zh-cn) and F (de-de), rows 2-4
#      - zh-cn sheet:     column C (Status), rows 2-4
#      - de-de sheet:     column C (Status), rows 2-4
#
# 2) Column width update for those same Status/locale columns:
#    17.2159881591797 -> 13.4101845877511 (narrower column)
#      - Overview sheet: columns E and F
#      - zh-cn sheet:     column C
#      - de-de sheet:     column C

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("E4").Value = $newStatus
$wsOverview.Range("F4").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("C4").Value = $newStatus

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("C4").Value = $newStatus

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
